$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 39 becomes the "Unassigned" row (values from old row 40, text capitalized)
$ws.Cells.Item(39, 1).Value = "Unassigned"
$ws.Cells.Item(39, 2).Value = "Unassigned"
$ws.Cells.Item(39, 3).Value = "Unassigned"
$ws.Cells.Item(39, 4).Value = 0
$ws.Cells.Item(39, 5).Value = 298
$ws.Cells.Item(39, 6).Value = 239

# Row 40 becomes the "Urophycis sp" row (values from old row 39)
$ws.Cells.Item(40, 1).Value = "Urophycis sp"
$ws.Cells.Item(40, 2).Value = "Red White or Spotted hake"
$ws.Cells.Item(40, 3).Value = "Teleost Fish"
$ws.Cells.Item(40, 4).Value = 1261
$ws.Cells.Item(40, 5).Value = 1025
$ws.Cells.Item(40, 6).Value = 11
